$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Centraal Station" stop to "Centraal Station_B"
$ws.Range("A18").Value = "Centraal Station_B"

# Move the active selection to A19 (just below the edited cell),
# matching the post-edit cursor position seen after pressing Enter.
$ws.Range("A19").Select()
